$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": row 59 - Registration run for cert89 gets re-styled to match
# the live-run formatting (column default) and the run-time value is refined
# to the more precise timestamp.
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# A59 holds a literal date-like string ("2022-09-06") - write it through a
# text formula and paste back as a value so it is not auto-converted into a
# date serial number.
$wsAmsin.Range("A59").Formula = "=""2022-09-06"""
$wsAmsin.Range("A59").Copy()
$wsAmsin.Range("A59").PasteSpecial(-4163)
$wsAmsin.Range("A59").Style = "Normal"
$wsAmsin.Range("A59").VerticalAlignment = -4107

# B59 - refine the recorded run time value.
$wsAmsin.Range("B59").Value = 44810.93921685185

# C59:G59 - re-apply the standard (column default) formatting.
$wsAmsin.Range("C59").Value = "cert89"
$wsAmsin.Range("D59").Value = 51
$wsAmsin.Range("E59").Value = 51
$wsAmsin.Range("F59").Value = 0
$wsAmsin.Range("G59").Value = 1.15
$wsAmsin.Range("C59:G59").Style = "Normal"
$wsAmsin.Range("C59:G59").VerticalAlignment = -4107

# ---------------------------------------------------------------------------
# Sheet "AMS": append the two new live registration-script runs (rows 16 & 17)
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# --- Row 16 : cert166 --------------------------------------------------
$wsAms.Range("A16").Formula = "=""2022-09-08"""
$wsAms.Range("A16").Copy()
$wsAms.Range("A16").PasteSpecial(-4163)

$wsAms.Range("B16").Value = 44812.4908575463
$wsAms.Range("B16").Style = "Normal"
$wsAms.Range("B16").NumberFormat = "yyyy-mm-dd h:mm:ss"
$wsAms.Range("B16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wsAms.Range("C16").Value = "cert166"
$wsAms.Range("D16").Value = 51
$wsAms.Range("E16").Value = 49
$wsAms.Range("F16").Value = 2
$wsAms.Range("G16").Value = 1.55

# --- Row 17 : certi166 --------------------------------------------------
$wsAms.Range("A17").Formula = "=""2022-09-08"""
$wsAms.Range("A17").Copy()
$wsAms.Range("A17").PasteSpecial(-4163)

$wsAms.Range("B17").Value = 44812.5083360391
$wsAms.Range("B17").Style = "Normal"
$wsAms.Range("B17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wsAms.Range("C17").Value = "certi166"
$wsAms.Range("D17").Value = 51
$wsAms.Range("E17").Value = 51
$wsAms.Range("F17").Value = 0
$wsAms.Range("G17").Value = 0.92
